$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 2770.8572
$ws.Range("I32").Value = 2899.5
$ws.Range("K32").Value = 2899.5
$ws.Range("M32").Value = -2573.5

$ws.Range("H40").Value = 5567.3335
$ws.Range("I40").Value = 3729.5833
$ws.Range("K40").Value = 3729.5833
$ws.Range("M40").Value = -3554.5833

$ws.Range("H58").Value = 2942.3333
$ws.Range("J58").Value = 4993
$ws.Range("L58").Value = 14979
$ws.Range("N58").Value = -15279

$ws.Range("H68").Value = 49998
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").Value = $null

$ws.Range("H71").Value = 49998
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").Value = $null

$ws.Range("H132").Value = 2547.4255
$ws.Range("I132").Value = 2547.4255
$ws.Range("K132").Value = 7642.2765
$ws.Range("M132").Value = -5112.2765

$ws.Range("H137").Value = 3574.0679
$ws.Range("I137").Value = 3093.8723
$ws.Range("J137").Value = 5454.8335
$ws.Range("K137").Value = 9281.616900000001
$ws.Range("L137").Value = 16364.5005
$ws.Range("M137").Value = -6731.616900000001
$ws.Range("N137").Value = -21464.5005

$ws.Range("H138").Value = 6047.59
$ws.Range("I138").Value = 6202.55
$ws.Range("J138").Value = 5994.1553
$ws.Range("K138").Value = 18607.65
$ws.Range("L138").Value = 17982.4659
$ws.Range("M138").Value = -13467.65
$ws.Range("N138").Value = -28262.4659

$ws.Range("H141").Value = 2241.1667
$ws.Range("I141").Value = 1719.4
$ws.Range("K141").Value = 5158.200000000001
$ws.Range("M141").Value = 21.79999999999927

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 3429.3333
$ws.Range("I31").Value = 3429.3333
$ws.Range("K31").Value = 3429.3333
$ws.Range("M31").Value = -3135.3333

$ws.Range("H32").Value = 22990.38
$ws.Range("I32").Value = 17012.139
$ws.Range("J32").Value = 59713.855
$ws.Range("K32").Value = 17012.139
$ws.Range("L32").Value = 59713.855
$ws.Range("M32").Value = -16725.139
$ws.Range("N32").Value = -60287.855

$ws.Range("H74").Value = 4464.1816
$ws.Range("I74").Value = 3571.3076
$ws.Range("J74").Value = 7780.5713
$ws.Range("K74").Value = 3571.3076
$ws.Range("L74").Value = 7780.5713
$ws.Range("M74").Value = -2697.3076
$ws.Range("N74").Value = -9528.5713

$ws.Range("H77").Value = 4464.1816
$ws.Range("I77").Value = 3571.3076
$ws.Range("J77").Value = 7780.5713
$ws.Range("K77").Value = 17856.538
$ws.Range("L77").Value = 38902.85649999999
$ws.Range("M77").Value = -13488.538
$ws.Range("N77").Value = -47638.85649999999

$ws.Range("H122").Value = 10397.714
$ws.Range("I122").Value = 9073.556
$ws.Range("K122").Value = 27220.668
$ws.Range("M122").Value = -24770.668

$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").Value = $null

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H36").Value = 4681.6
$ws.Range("I36").Value = 852
$ws.Range("K36").Value = 852
$ws.Range("M36").Value = -318

$ws.Range("H86").Value = 18581
$ws.Range("I86").Value = 18900.5
$ws.Range("J86").Value = 18197.6
$ws.Range("K86").Value = 18900.5
$ws.Range("L86").Value = 18197.6
$ws.Range("M86").Value = -17777.5
$ws.Range("N86").Value = -20443.6

$ws.Range("H89").Value = 18581
$ws.Range("I89").Value = 18900.5
$ws.Range("J89").Value = 18197.6
$ws.Range("K89").Value = 94502.5
$ws.Range("L89").Value = 90988
$ws.Range("M89").Value = -88886.5
$ws.Range("N89").Value = -102220

$ws.Range("H102").Value = 7499.5
$ws.Range("I102").Value = 7499.5
$ws.Range("K102").Value = 7499.5
$ws.Range("M102").Value = -4254.5

$ws.Range("H132").Value = 82666.336
$ws.Range("J132").Value = 82666.336
$ws.Range("L132").Value = 82666.336
$ws.Range("N132").Value = -92786.336

$ws.Range("H134").Value = 6160.3423
$ws.Range("I134").Value = 5161.4814
$ws.Range("J134").Value = 8612.091
$ws.Range("K134").Value = 15484.4442
$ws.Range("L134").Value = 25836.273
$ws.Range("M134").Value = -12949.4442
$ws.Range("N134").Value = -30906.273

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7858.326
$ws.Range("I31").Value = 4359.8125
$ws.Range("J31").Value = 15854.929
$ws.Range("K31").Value = 4359.8125
$ws.Range("L31").Value = 15854.929
$ws.Range("M31").Value = -4064.8125
$ws.Range("N31").Value = -16444.929

$ws.Range("H34").Value = 7858.326
$ws.Range("I34").Value = 4359.8125
$ws.Range("J34").Value = 15854.929
$ws.Range("K34").Value = 4359.8125
$ws.Range("L34").Value = 15854.929
$ws.Range("M34").Value = -4157.8125
$ws.Range("N34").Value = -16258.929

$ws.Range("H51").Value = 51999
$ws.Range("J51").Value = 51999
$ws.Range("L51").Value = 51999
$ws.Range("N51").Value = -53471

$ws.Range("H59").Value = 56666.332
$ws.Range("J59").Value = 79999.5
$ws.Range("L59").Value = 79999.5
$ws.Range("N59").Value = -82289.5

$ws.Range("H60").Value = 44177.332
$ws.Range("I60").Value = 42101
$ws.Range("J60").Value = 44436.875
$ws.Range("K60").Value = 42101
$ws.Range("L60").Value = 44436.875
$ws.Range("M60").Value = -41590
$ws.Range("N60").Value = -45458.875

$ws.Range("H61").Value = 51999
$ws.Range("J61").Value = 51999
$ws.Range("L61").Value = 51999
$ws.Range("N61").Value = -52695

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 960.8
$ws.Range("I8").Value = 960.8
$ws.Range("K8").Value = 2882.4
$ws.Range("M8").Value = -2743.4

$ws.Range("H23").Value = 214.2
$ws.Range("J23").Value = 181.14285
$ws.Range("L23").Value = 543.4285500000001
$ws.Range("N23").Value = -1013.42855

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3341.9473
$ws.Range("I102").Value = 3009.8572
$ws.Range("K102").Value = 3009.8572
$ws.Range("M102").Value = -1387.8572

$ws.Range("H132").Value = 7777.604
$ws.Range("I132").Value = 6985.6216
$ws.Range("J132").Value = 10441.546
$ws.Range("K132").Value = 20956.8648
$ws.Range("L132").Value = 31324.638
$ws.Range("M132").Value = -18426.8648
$ws.Range("N132").Value = -36384.638

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1174.2963
$ws.Range("I16").Value = 1017.8
$ws.Range("K16").Value = 1017.8
$ws.Range("M16").Value = -847.8

$ws.Range("H122").Value = 9975
$ws.Range("I122").Value = 5900
$ws.Range("K122").Value = 17700
$ws.Range("M122").Value = -15250

$ws.Range("H132").Value = 3698.2812
$ws.Range("I132").Value = 3698.2812
$ws.Range("K132").Value = 11094.8436
$ws.Range("M132").Value = -8564.8436

$ws.Range("H136").Value = 9144.566000000001
$ws.Range("I136").Value = 6223.8335
$ws.Range("J136").Value = 9874.75
$ws.Range("K136").Value = 18671.5005
$ws.Range("L136").Value = 29624.25
$ws.Range("M136").Value = -16121.5005
$ws.Range("N136").Value = -34724.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 48039.125
$ws.Range("J52").Value = 42999
$ws.Range("L52").Value = 42999
$ws.Range("N52").Value = -43451

$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").Value = $null

$ws.Range("H81").Value = 400.5
$ws.Range("I81").Value = 400.5
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 801
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = 260
$ws.Range("N81").Value = $null

$ws.Range("H84").Value = 400.5
$ws.Range("I84").Value = 400.5
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 4005
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = 1299
$ws.Range("N84").Value = $null

$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").Value = $null

$ws.Range("H136").Value = 6206.143
$ws.Range("I136").Value = 5082.75
$ws.Range("K136").Value = 15248.25
$ws.Range("M136").Value = -12698.25

